$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.368.90"
$ws.Range("E2").Value = "  -3.91%  "

$ws.Range("D3").Value = "3.570.26"
$ws.Range("E3").Value = "  -4.27%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.84%  "

$ws.Range("D7").Value = "3.563.69"
$ws.Range("E7").Value = "  -4.42%  "

$ws.Range("E8").Value = "  -4.18%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.670"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.08%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.147"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.85%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.99%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000257"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -12.71%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -8.47%  "

$ws.Range("D15").Value = "4.139.73"
$ws.Range("E15").Value = "  -4.24%  "

$ws.Range("D16").Value = "3.565.10"
$ws.Range("E16").Value = "  -4.53%  "

$ws.Range("E17").Value = "  -0.97%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.51%  "

$ws.Range("E19").Value = "  -5.98%  "

$ws.Range("D20").Value = "66.198.89"
$ws.Range("E20").Value = "  -3.92%  "

$ws.Range("E21").Value = "  -7.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "396.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.21%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.07%  "

$ws.Range("E26").Value = "  -5.99%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.13%  "

$ws.Range("E28").Value = "  -0.72%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.68%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.90%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.95%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.99%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "616.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.72%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "63.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.23%  "

$ws.Range("E36").Value = "  -8.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "41.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.20%  "

$ws.Range("E38").Value = "  +0.26%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.392"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.16%  "

$ws.Range("D40").Value = "0.0₃0753"
$ws.Range("E40").Value = "  -13.85%  "

$ws.Range("E41").Value = "  -7.89%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.06%  "

$ws.Range("D43").Value = "2.999.47"
$ws.Range("E43").Value = "  +5.93%  "

$ws.Range("E44").Value = "  -8.52%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.48%  "

$ws.Range("E46").Value = "  -8.58%  "

$ws.Range("E47").Value = "  -6.87%  "

$ws.Range("E48").Value = "  -1.73%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.01%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.17%  "

$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "138.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.01%  "

